$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.843.52'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '1.763.38'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '328.26'
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  -1.71%  '
$ws.Range("D8").Value = '0.3544'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '0.07402'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '41.89'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = '1.099'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").Value = '7.231'
$ws.Range("E15").Value = '  +2.83%  '
$ws.Range("D16").Value = '1.758.83'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '92.88'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '0.06427'
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '17.05'
$ws.Range("E21").Value = '  +3.27%  '
$ws.Range("D22").Value = '5.759'
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").Value = '27.875.50'
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("D25").Value = '2.108'
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").Value = '160.34'
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("D28").Value = '1.963.31'
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '2.140'
$ws.Range("E29").Value = '  +5.45%  '
$ws.Range("D30").Value = '123.97'
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("D31").Value = '1.099'
$ws.Range("E31").Value = '  +5.47%  '
$ws.Range("D32").Value = '0.09188'
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").Value = '5.638'
$ws.Range("E33").Value = '  +5.20%  '
$ws.Range("D34").Value = '3.690'
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").Value = '11.84'
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("D36").Value = '0.06176'
$ws.Range("E36").Value = '  +4.02%  '
$ws.Range("D37").Value = '0.02276'
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("D38").Value = '0.2097'
$ws.Range("E38").Value = '  +2.27%  '
$ws.Range("D39").Value = '0.6305'
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("D40").Value = '4.945'
$ws.Range("E40").Value = '  +1.64%  '
$ws.Range("D41").Value = '1.182'
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").Value = '1.391'
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").Value = '7.831'
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("D44").Value = '13.23'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").Value = '3.746'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").Value = '0.5838'
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("D47").Value = '122.19'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = '1.949'
$ws.Range("D49").Value = '0.06877'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").Value = '72.66'
$ws.Range("E51").Value = '  +2.37%  '
